# Applies the "Add/update resource data for BiblicaStudyNotesKeyTerms" edit.
$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1) Remove the whole "License Information" heading paragraph.
# ----------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("License Information", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Expand(4)  # wdParagraph - include the paragraph mark
    $rng.Delete()
}

# ----------------------------------------------------------------------
# 2) Rewrite the "關鍵詞 (Biblica) ... license." paragraph into the new
#    Biblica Study Notes (Key Terms) description (also removes the two
#    hyperlinks that used to be embedded in the old text).
# ----------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("is based on", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Expand(4)  # wdParagraph
    $pStart = $rng.Start
    $pEnd = $rng.End
    $scoped = $d.Range($pStart, $pEnd)

    $oldAll = "關鍵詞 (Biblica) (Chinese (Traditional)) is based on: Biblica Bible Dictionary, Biblica, Inc., 2023, which is licensed under a CC BY-SA 4.0 license."
    $boldPart = "Biblica Study Notes (Key Terms)"
    $restPart = " © 2023 Biblica Inc. Released under CC BY-SA 4.0 license. Biblica Study Notes has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文)from Biblica Study Notes © 2023 Biblica Inc. Released under CC BY-SA 4.0 license by Mission Mutual."
    $newAll = $boldPart + $restPart

    $replaced = $scoped.Find.Execute($oldAll, $false, $false, $false, $false, $false, $true, 1, $false, $newAll, 2)

    if ($replaced) {
        # The whole replacement text lands inside the old bold run, so the
        # text that must stay plain needs to be un-bolded explicitly.
        $boldEnd = $scoped.Start + $boldPart.Length
        $plainRange = $d.Range($boldEnd, $scoped.End)
        $plainRange.Font.Bold = 0
    }
}

# ----------------------------------------------------------------------
# 3) Remove the whole "This PDF version is provided under the same
#    license." paragraph.
# ----------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("This PDF version is provided under the same license.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Expand(4)  # wdParagraph
    $rng.Delete()
}

# ----------------------------------------------------------------------
# 4) Remove the whole italic "幻影說" paragraph that used to sit right
#    after the "huan" heading.
# ----------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("huan", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $afterHuan = $d.Range($rng.End, $d.Content.End)
    $found2 = $afterHuan.Find.Execute("幻影說", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found2) {
        $afterHuan.Expand(4)  # wdParagraph
        $afterHuan.Delete()
    }
}
